$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C7-Result row (row 12): update Pages to 10 and Comment to "CHECKING LATER"
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = "CHECKING LATER"

# Recalculate so the SUM formula in C13 updates
$excel.Calculate()

# Update the active selection to C13, matching the saved workbook view
$ws.Range("C13").Select()
